$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. "42.995.75" using
# "." as a thousands separator). Force text format first so Excel does not
# reinterpret/round these strings as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.995.75"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.300.98"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.18"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.95"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.517"
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.01"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.117"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.80"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.660.05"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.288.83"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.779"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.954.26"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.63"
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.11"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.21"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.05"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.21"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.79"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.08"
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.07"
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.01"
$ws.Range("E34").Value = "  -3.54%  "
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.69"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0686"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.110"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.005.03"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.17"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.46"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.79"
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.90"
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.58"
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.527.01"
$ws.Range("E51").Value = "  -0.81%  "
